# Add a new "2021年" row (row 5) of data below the existing 2018-2020 rows
# on Sheet1, following the same layout as the existing rows (A = year label,
# B:DK = the per-industry growth figures). Column AG has no data for 2021
# (same as it is blank for the prior years), so it is intentionally left
# untouched/empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "2021年"

$values = @(
    6.7, 28, -21.7, 40.3, 3, 6.3, 13.6, 8.5, 8.5, -7.7,
    -13.3, 8.699999999999999, -41.1, -2.3, 10.1, -13.2, 63.1, 12.3, 14.8, 11.5,
    4.9, 24.2, 15.1, 19.7, 83.8, 17.2, 21, 27.6, -1.9, 13.7,
    -44.3, $null, -23.4, -13.4, -9.300000000000001, 18.5, 3.9, 2, -11.3, -11.5,
    -26.1, 11.5, -4.9, 46.3, 28.3, 24.6, 10.3, 6.4, -4.9, -4.7,
    11.9, 1.6, 6.4, 10.9, 48.6, 2.5, 8.699999999999999, 21.6, -7.5, -1.8,
    13.5, 19.1, -1.5, 3.8, 1.1, -6, 7.8, -36.3, 9.6, -1.5,
    -0.4, -10.6, 4.9, 6.6, 22.9, 21, -2.7, 5.8, 59.2, 13.4,
    -30.1, -3.2, 13.4, 17, 81.59999999999999, 14.6, 0.7, 18.5, -3.5, -17,
    1.1, 12.8, 21.1, 4.8, -25, 12.7, 7.9, 26.1, 0.6, 44.5,
    15.8, 36.9, 33.9, 10.5, -2.2, 24.5, 7.2, -4.6, 13.6, 38.7,
    8.1, 7.8, 4.5, 27.9
)

# Values start in column B (2) and run through column DK.
for ($i = 0; $i -lt $values.Count; $i++) {
    $col = $i + 2
    if ($values[$i] -ne $null) {
        $ws.Cells.Item(5, $col).Value = $values[$i]
    }
}
